$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting of the previous row down to the new row, then fill in values
$ws.Range("A3:F3").Copy()
$ws.Range("A4:F4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "?"
$ws.Range("C4").Value = "Plataforma para venda e compra de ingressos"
$ws.Range("D4").Value = "Vinícius Radé"
$ws.Range("E4").Value2 = 42544
$ws.Range("F4").Value = "Cadastrado"
